# "clean up of prob programming"
# Update the lesson-log worksheet: rewrite the MCMC/prob-programming entries
# in column F (rows 15-16), add three new weeks of content (rows 17-20,
# including the E-column date formulas), and refresh the sheet view
# (scroll position, zoom, selection) to match where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: shorten/replace the F-column note; its wrap height shrinks
#     from 3 lines to 2 because the long text moved out of this cell.
$ws.Range("F15").Value = " MCMC facendo poi implementare island hopping"
$ws.Rows.Item(15).RowHeight = 34

# --- Row 18's new cell reuses row 16's current highlight-fill format, so
#     grab a copy of it before row 16 gets its own formatting reset below.
$ws.Range("F16").Copy()
$ws.Range("F18").PasteSpecial(-4122)

# --- Row 16: new note text, and the old highlight fill is cleared
#     (copy the plain formatting from F14, a same-column cell with the
#     default/no style, then set the value).
$ws.Range("F14").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16").Value = "chiudere MCMC, prob programming"

# --- Row 17: new F entry, formatted like F13/F9 ("Neutral" highlight style).
$ws.Range("F13").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value = "exe MF"

# --- Row 18: fill in the E-column date formula (continues the +7 weekly
#     series) and set the F-column note text (format already copied above).
$ws.Range("E18").Formula = "=E16+7"
$ws.Range("F18").Value = "prob programming, presentare  progetto"

# --- Row 19: E-column date formula + plain F-column note (default style).
$ws.Range("E19").Formula = "=E17+7"
$ws.Range("F19").Value = "normal-normal"

# --- Row 20: E-column date formula + plain F-column note (default style).
$ws.Range("E20").Formula = "=E18+7"
$ws.Range("F20").Value = "hyp test"

# --- Sheet view: scroll position, zoom, and selection as left by the author.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$win.Zoom = 125
$ws.Range("F21").Select()
